$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3350.5
$ws.Range("I29").Value = 103
$ws.Range("J29").Value = 4000
$ws.Range("K29").Value = 309
$ws.Range("L29").Value = 12000
$ws.Range("M29").Value = -28
$ws.Range("N29").Value = -12562

$ws.Range("H33").Value = 461.85715
$ws.Range("I33").Value = 105.947365
$ws.Range("J33").Value = 1213.2222
$ws.Range("K33").Value = 105.947365
$ws.Range("L33").Value = 1213.2222
$ws.Range("M33").Value = 123.052635
$ws.Range("N33").Value = -1671.2222

$ws.Range("H121").Value = 1235.8636
$ws.Range("J121").Value = 1289.45
$ws.Range("L121").Value = 3868.35
$ws.Range("N121").Value = -7362.35

$ws.Range("H132").Value = 1636462.2
$ws.Range("I132").Value = 3097.48
$ws.Range("J132").Value = 9803286
$ws.Range("K132").Value = 9292.440000000001
$ws.Range("L132").Value = 29409858
$ws.Range("M132").Value = -6762.440000000001
$ws.Range("N132").Value = -29414918

$ws.Range("H135").Value = 19298.893
$ws.Range("I135").Value = 24306.883
$ws.Range("J135").Value = 2734
$ws.Range("K135").Value = 218761.947
$ws.Range("L135").Value = 24606
$ws.Range("M135").Value = -216226.947
$ws.Range("N135").Value = -29676

$ws.Range("H137").Value = 1924462.5
$ws.Range("I137").Value = 2858226
$ws.Range("J137").Value = 2008.2354
$ws.Range("K137").Value = 8574678
$ws.Range("L137").Value = 6024.706200000001
$ws.Range("M137").Value = -8572128
$ws.Range("N137").Value = -11124.7062

$ws.Range("H140").Value = 67778.17999999999
$ws.Range("J140").Value = 67778.17999999999
$ws.Range("L140").Value = 67778.17999999999
$ws.Range("N140").Value = -78138.17999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 15468.777
$ws.Range("J42").Value = 15468.777
$ws.Range("L42").Value = 15468.777
$ws.Range("N42").Value = -16440.777

$ws.Range("H61").Value = 20041002
$ws.Range("I61").Value = 21761840
$ws.Range("J61").Value = 251353.5
$ws.Range("K61").Value = 21761840
$ws.Range("L61").Value = 251353.5
$ws.Range("M61").Value = -21761628
$ws.Range("N61").Value = -251777.5

$ws.Range("I74").Value = 7607267
$ws.Range("J74").Value = 92927.17999999999
$ws.Range("K74").Value = 7607267
$ws.Range("L74").Value = 92927.17999999999
$ws.Range("M74").Value = -7606393
$ws.Range("N74").Value = -94675.17999999999

$ws.Range("I77").Value = 7607267
$ws.Range("J77").Value = 92927.17999999999
$ws.Range("K77").Value = 38036335
$ws.Range("L77").Value = 464635.9
$ws.Range("M77").Value = -38031967
$ws.Range("N77").Value = -473371.9

$ws.Range("H102").Value = 6809254.5
$ws.Range("I102").Value = 7525734
$ws.Range("J102").Value = 2700
$ws.Range("K102").Value = 7525734
$ws.Range("L102").Value = 2700
$ws.Range("M102").Value = -7524112
$ws.Range("N102").Value = -5944

$ws.Range("H132").Value = 51770.195
$ws.Range("I132").Value = 28552.445
$ws.Range("K132").Value = 85657.33499999999
$ws.Range("M132").Value = -83127.33499999999

$ws.Range("H135").Value = 36861.285
$ws.Range("J135").Value = 36861.285
$ws.Range("L135").Value = 36861.285
$ws.Range("N135").Value = -47001.285

$ws.Range("H136").Value = 20041002
$ws.Range("I136").Value = 21761840
$ws.Range("J136").Value = 251353.5
$ws.Range("K136").Value = 65285520
$ws.Range("L136").Value = 754060.5
$ws.Range("M136").Value = -65282970
$ws.Range("N136").Value = -759160.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14251.947
$ws.Range("I86").Value = 17304.867
$ws.Range("J86").Value = 2803.5
$ws.Range("K86").Value = 17304.867
$ws.Range("L86").Value = 2803.5
$ws.Range("M86").Value = -16181.867
$ws.Range("N86").Value = -5049.5

$ws.Range("H89").Value = 14251.947
$ws.Range("I89").Value = 17304.867
$ws.Range("J89").Value = 2803.5
$ws.Range("K89").Value = 86524.33499999999
$ws.Range("L89").Value = 14017.5
$ws.Range("M89").Value = -80908.33499999999
$ws.Range("N89").Value = -25249.5

$ws.Range("H94").Value = 1015.5714
$ws.Range("I94").Value = 527.25
$ws.Range("J94").Value = 1666.6666
$ws.Range("K94").Value = 527.25
$ws.Range("L94").Value = 1666.6666
$ws.Range("M94").Value = -76.25
$ws.Range("N94").Value = -2568.6666

$ws.Range("H99").Value = 795.8823
$ws.Range("I99").Value = 679.2308
$ws.Range("K99").Value = 679.2308
$ws.Range("M99").Value = 818.7692

$ws.Range("H105").Value = 50001980
$ws.Range("I105").Value = 166668200
$ws.Range("J105").Value = 2171.4285
$ws.Range("K105").Value = 166668200
$ws.Range("L105").Value = 2171.4285
$ws.Range("M105").Value = -166666453
$ws.Range("N105").Value = -5665.4285

$ws.Range("H107").Value = 1452.8572
$ws.Range("I107").Value = 1151.8889
$ws.Range("K107").Value = 1151.8889
$ws.Range("M107").Value = 768.1111000000001

$ws.Range("H134").Value = 1459.8276
$ws.Range("I134").Value = 822.55316
$ws.Range("J134").Value = 4182.727
$ws.Range("K134").Value = 2467.65948
$ws.Range("L134").Value = 12548.181
$ws.Range("M134").Value = 67.34051999999974
$ws.Range("N134").Value = -17618.181

$ws.Range("H140").Value = 55100
$ws.Range("J140").Value = 55100
$ws.Range("L140").Value = 55100
$ws.Range("N140").Value = -65460

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1977.625
$ws.Range("I105").Value = 1968.3334
$ws.Range("J105").Value = 2005.5
$ws.Range("K105").Value = 1968.3334
$ws.Range("L105").Value = 2005.5
$ws.Range("M105").Value = -221.3334
$ws.Range("N105").Value = -5499.5

$ws.Range("H107").Value = 322.24
$ws.Range("I107").Value = 327.90475
$ws.Range("J107").Value = 292.5
$ws.Range("K107").Value = 327.90475
$ws.Range("L107").Value = 292.5
$ws.Range("M107").Value = 1592.09525
$ws.Range("N107").Value = -4132.5

$ws.Range("H132").Value = 156989.39
$ws.Range("I132").Value = 114673.555
$ws.Range("J132").Value = 252200
$ws.Range("K132").Value = 344020.665
$ws.Range("L132").Value = 756600
$ws.Range("M132").Value = -341490.665
$ws.Range("N132").Value = -761660

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 329
$ws.Range("I6").Value = 83.833336
$ws.Range("K6").Value = 251.500008
$ws.Range("M6").Value = -138.500008

$ws.Range("H131").Value = 1148.6522
$ws.Range("J131").Value = 1179.7273
$ws.Range("L131").Value = 3539.1819
$ws.Range("N131").Value = -13619.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5687500
$ws.Range("I11").Value = 6428571.5
$ws.Range("K11").Value = 6428571.5
$ws.Range("M11").Value = -6428432.5

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H13").Value = 523
$ws.Range("I13").Value = 523
$ws.Range("K13").Value = 523
$ws.Range("M13").Value = -384

$ws.Range("H132").Value = 41222.98
$ws.Range("I132").Value = 24988.928
$ws.Range("J132").Value = 126451.75
$ws.Range("K132").Value = 74966.784
$ws.Range("L132").Value = 379355.25
$ws.Range("M132").Value = -72436.784
$ws.Range("N132").Value = -384415.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2383.3333
$ws.Range("J2").Value = 2800
$ws.Range("L2").Value = 2800
$ws.Range("N2").Value = -3024

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H122").Value = 2957.4119
$ws.Range("I122").Value = 2291.3076
$ws.Range("J122").Value = 3369.762
$ws.Range("K122").Value = 6873.9228
$ws.Range("L122").Value = 10109.286
$ws.Range("M122").Value = -4423.9228
$ws.Range("N122").Value = -15009.286

$ws.Range("H132").Value = 33386.953
$ws.Range("I132").Value = 22658.307
$ws.Range("J132").Value = 66243.44
$ws.Range("K132").Value = 67974.921
$ws.Range("L132").Value = 198730.32
$ws.Range("M132").Value = -65444.921
$ws.Range("N132").Value = -203790.32

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3308.4211
$ws.Range("I122").Value = 2951.2222
$ws.Range("J122").Value = 3629.9
$ws.Range("K122").Value = 8853.6666
$ws.Range("L122").Value = 10889.7
$ws.Range("M122").Value = -6403.6666
$ws.Range("N122").Value = -15789.7

$ws.Range("H129").Value = 41741.4
$ws.Range("J129").Value = 41741.4
$ws.Range("L129").Value = 41741.4
$ws.Range("N129").Value = -51741.4

$ws.Range("H132").Value = 39869.277
$ws.Range("I132").Value = 32866.773
$ws.Range("K132").Value = 98600.319
$ws.Range("M132").Value = -96070.319
